$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 71, shifting existing rows 71-135 down to 72-136.
$ws.Rows.Item(71).Insert()

# Populate the new row 71 with its data.
$ws.Range("A71").Value = 4
$ws.Range("B71").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value = "Los Lagos"
$ws.Range("D71").Value = 44781
$ws.Range("D71").NumberFormat = $ws.Range("D72").NumberFormat
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = 100112052
$ws.Range("G71").Value = "Albahaca"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 70
$ws.Range("K71").Value = 6500
$ws.Range("L71").Value = 6500
$ws.Range("M71").Value = 6500
$ws.Range("N71").Value = '$/paquete'
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 6500
$ws.Range("Q71").Value = 1
$ws.Range("R71").Value = "Hortaliza"
